# #SSW-1990 update Standard Importe
#
# The import template's "Geburtsdatum" (birth date) sample column was being
# read back in as a real Excel date serial, which broke downstream import
# parsing for some locales. Switch the whole sample row to plain Text
# formatting and store the sample date as literal text instead of a date
# value, then leave the selection where the user last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Force the header row and the sample data row onto Text number format
# ("@") instead of the previous mix of General / date (mm-dd-yy) formats.
$ws.Range("A1:F1").NumberFormat = "@"
$ws.Range("A2:F2").NumberFormat = "@"

# The sample "Geburtsdatum" cell used to hold a real date value
# (01.01.2020 as a date serial); store it as literal text now so imports
# that expect a fixed-width text column aren't tripped up by regional
# date formatting.
$ws.Range("C2").Value = "01.01.2020"

# Restore the cursor/selection to where it was left (D10) when the file
# was last saved.
[void]$ws.Range("D10").Select()
